# Update the "ProjectSchedule" Gantt chart worksheet: replace the five
# placeholder tasks ("Task 1".."Task 5") with the actual software-phase
# tasks for the project, and correct a few start/end dates to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# Rename the placeholder tasks to the real task names.
$ws.Range("B19").Value = "Learn Arduino"
$ws.Range("B20").Value = "Block detection software"
$ws.Range("B21").Value = "Block sorting software"
$ws.Range("B22").Value = "Adding redundancies and updating code"
$ws.Range("B23").Value = "Testing feasibility of different means of detection"

# Overwrite a few date cells with fixed (literal) dates, replacing the
# formulas that used to compute them. Values are Excel date serials.
$ws.Range("F19").Value = 44872   # 07/11/2022 - end date for "Learn Arduino"
$ws.Range("F20").Value = 44879   # 14/11/2022 - end date for "Block detection software"
$ws.Range("E23").Value = 44871   # 06/11/2022 - start date for "Testing feasibility..."

$excel.Calculate()
